# Weekly update: insert a new week's record for
# "Comercializadora del Agro de Limarí - Zapallo italiano" as the new
# row 6, pushing the existing historical rows (old rows 6-38) down by
# one row (new rows 7-39).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 6; Excel shifts rows 6:38 -> 7:39
# and carries the row formatting (e.g. the date style on column D) down
# with them.
$ws.Rows("6:6").Insert()

# Populate the newly inserted row 6 with this week's data.
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C6").Value = "Coquimbo"
$ws.Range("D6").Value = 44532
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 100112032
$ws.Range("G6").Value = "Zapallo italiano"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 600
$ws.Range("K6").Value = 5500
$ws.Range("L6").Value = 6000
$ws.Range("M6").Value = 5750
$ws.Range("N6").Value = "$/caja 60 unidades"
$ws.Range("O6").Value = "Provincia de Limarí"
$ws.Range("P6").Value = 96
$ws.Range("Q6").Value = 60
$ws.Range("R6").Value = "Hortaliza"
